$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates to the Price (D) and Volume(1h) (E) columns.
$updates = @(
    @{ Row = 2; D = '62.843.83'; E = '  -0.52%  ' }
    @{ Row = 3; D = '3.220.24'; E = '  -0.24%  ' }
    @{ Row = 4; D = $null; E = '  -0.85%  ' }
    @{ Row = 5; D = '526.83'; E = '  +5.28%  ' }
    @{ Row = 6; D = '171.90'; E = '  -0.38%  ' }
    @{ Row = 7; D = '0.594'; E = '  +3.22%  ' }
    @{ Row = 8; D = $null; E = '  -0.25%  ' }
    @{ Row = 9; D = '3.215.51'; E = '  -0.19%  ' }
    @{ Row = 10; D = $null; E = '  +1.29%  ' }
    @{ Row = 11; D = '52.95'; E = '  -4.76%  ' }
    @{ Row = 12; D = $null; E = '  +5.38%  ' }
    @{ Row = 13; D = '0.0000251'; E = '  +2.43%  ' }
    @{ Row = 14; D = '9.09'; E = '  +3.86%  ' }
    @{ Row = 15; D = '3.736.72'; E = '  -0.62%  ' }
    @{ Row = 16; D = $null; E = '  -4.04%  ' }
    @{ Row = 17; D = '3.213.91'; E = '  -1.02%  ' }
    @{ Row = 18; D = '62.820.91'; E = '  -0.60%  ' }
    @{ Row = 19; D = '17.14'; E = '  +2.92%  ' }
    @{ Row = 20; D = '10.99'; E = '  +4.84%  ' }
    @{ Row = 21; D = '0.966'; E = '  +5.50%  ' }
    @{ Row = 22; D = '365.11'; E = '  +1.66%  ' }
    @{ Row = 23; D = '3.75'; E = '  +5.96%  ' }
    @{ Row = 24; D = '81.17'; E = '  +4.04%  ' }
    @{ Row = 25; D = '10.98'; E = '  +4.76%  ' }
    @{ Row = 26; D = '3.92'; E = '  +7.63%  ' }
    @{ Row = 27; D = '6.13'; E = '  +2.10%  ' }
    @{ Row = 28; D = $null; E = '  +3.78%  ' }
    @{ Row = 29; D = '11.25'; E = '  +4.22%  ' }
    @{ Row = 30; D = '8.12'; E = '  +1.46%  ' }
    @{ Row = 31; D = '28.38'; E = '  +3.19%  ' }
    @{ Row = 32; D = '631.21'; E = '  -1.21%  ' }
    @{ Row = 33; D = '6.39'; E = '  +0.29%  ' }
    @{ Row = 34; D = '11.22'; E = '  +4.89%  ' }
    @{ Row = 35; D = '0.105'; E = '  +6.12%  ' }
    @{ Row = 36; D = '56.72'; E = '  -2.13%  ' }
    @{ Row = 37; D = $null; E = '  -0.04%  ' }
    @{ Row = 38; D = '36.59'; E = '  +6.17%  ' }
    @{ Row = 39; D = $null; E = '  +4.07%  ' }
    @{ Row = 40; D = $null; E = '  -0.60%  ' }
    @{ Row = 41; D = '0.0₃0703'; E = '  +16.29%  ' }
    @{ Row = 42; D = '0.122'; E = '  +3.53%  ' }
    @{ Row = 43; D = '2.867.95'; E = '  +3.82%  ' }
    @{ Row = 44; D = '2.53'; E = '  +15.73%  ' }
    @{ Row = 45; D = $null; E = '  +5.02%  ' }
    @{ Row = 46; D = $null; E = '  +15.59%  ' }
    @{ Row = 47; D = '0.0391'; E = '  +5.61%  ' }
    @{ Row = 48; D = $null; E = '  -1.20%  ' }
    @{ Row = 49; D = $null; E = '  +10.09%  ' }
    @{ Row = 50; D = '0.125'; E = '  +3.73%  ' }
    @{ Row = 51; D = '134.36'; E = '  +2.42%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D" + $u.Row)
        # Force text storage so numeric-looking strings (e.g. "526.83")
        # aren't silently reinterpreted as numbers by Excel, matching the
        # original inline-string cell content. Reset the style afterwards
        # so no stray formatting is introduced on the cell.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    $ws.Range("E" + $u.Row).Value = $u.E
}
